$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the monthly fee/charge values in column G for rows 18-20
# from 1423500 to 1300000 (parte 1 de nuevos estado de cuenta)
$ws.Range("G18").Value = 1300000
$ws.Range("G19").Value = 1300000
$ws.Range("G20").Value = 1300000

$wb.Save()
